# Updates cryptos list values (GitHub Actions style refresh of Price / Volume(1h)
# columns, plus a VeChain/Hedera row swap) to match the latest scrape.
#
# Note: several "Price" values look like plain numbers (e.g. "0.500", "593.59").
# Assigning those directly to .Value would make Excel auto-convert them to
# floating point numbers (losing exact text/trailing zeros). To keep them as
# literal text - exactly like the source inline strings - we prefix the value
# with a single quote (Excel's classic "force text" prefix) and then reset the
# cell's Style back to "Normal" so no stray Text-number-format style lingers
# on the cell (matching the original, un-styled price cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.940.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.64%  '
$ws.Range("D3").Value = '''3.213.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -8.37%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''593.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").Value = '''151.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -12.45%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''3.205.57'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -8.49%  '
$ws.Range("D9").Value = '''0.544'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -10.55%  '
$ws.Range("E10").Value = '  -11.06%  '
$ws.Range("D11").Value = '''6.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.29%  '
$ws.Range("D12").Value = '''0.500'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -14.52%  '
$ws.Range("D13").Value = '''38.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -15.54%  '
$ws.Range("E14").Value = '  -11.42%  '
$ws.Range("D15").Value = '''3.732.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.39%  '
$ws.Range("D16").Value = '''66.941.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.61%  '
$ws.Range("D17").Value = '''3.219.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.20%  '
$ws.Range("E18").Value = '  -4.98%  '
$ws.Range("D19").Value = '''7.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -13.73%  '
$ws.Range("D20").Value = '''530.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -13.69%  '
$ws.Range("D21").Value = '''14.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -14.24%  '
$ws.Range("D22").Value = '''0.760'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -13.48%  '
$ws.Range("D23").Value = '''7.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -13.15%  '
$ws.Range("D24").Value = '''13.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.24%  '
$ws.Range("D25").Value = '''85.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -13.09%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -14.60%  '
$ws.Range("D28").Value = '''2.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -14.11%  '
$ws.Range("D29").Value = '''8.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.10%  '
$ws.Range("D30").Value = '''29.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -13.52%  '
$ws.Range("D31").Value = '''2.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.04%  '
$ws.Range("D33").Value = '''544.88'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -13.93%  '
$ws.Range("D34").Value = '''6.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -18.95%  '
$ws.Range("D35").Value = '''5.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -16.05%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").Value = '''53.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.85%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.0426'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.87%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.0865'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.21%  '
$ws.Range("D40").Value = '''9.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -13.21%  '
$ws.Range("E41").Value = '  -12.13%  '
$ws.Range("D42").Value = '''2.919.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = '''2.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -24.26%  '
$ws.Range("D44").Value = '''0.263'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -15.22%  '
$ws.Range("D45").Value = '''0.0₃0583'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -20.61%  '
$ws.Range("E46").Value = '  -17.16%  '
$ws.Range("D47").Value = '''26.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -16.60%  '
$ws.Range("D49").Value = '''2.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -17.28%  '
$ws.Range("D50").Value = '''0.114'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -12.61%  '
$ws.Range("D51").Value = '''123.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.74%  '
